{"js": "// fix(docx): fix OOXMLValidator error on KeywordTok output\n//\n// The <w:rPr> child-element order inside several syntax-highlighting\n// character styles in styles.xml doesn't match the sequence required by\n// wml.xsd (CT_RPr) -- <w:b/>/<w:i/> must come before <w:color/>, e.g.:\n//   <w:rPr><w:color w:val=\"007020\"/><w:b/></w:rPr>   (wrong)\n//   <w:rPr><w:b/><w:color w:val=\"007020\"/></w:rPr>   (correct)\n// Re-assert each style's existing bold/italic flag so the host\n// re-serializes that style's <w:rPr> in canonical schema order (the\n// color itself is untouched).\n\n// styleId -> which flag(s) to re-assert. Only the flag(s) actually\n// present need touching; re-setting one forces the whole <w:rPr> to be\n// rewritten in schema order, fixing any other out-of-order sibling too.\nconst fixups = {\n  KeywordTok: { bold: true },\n  ImportTok: { bold: true },\n  CommentTok: { italic: true },\n  DocumentationTok: { italic: true },\n  AnnotationTok: { bold: true, italic: true },\n  CommentVarTok: { bold: true, italic: true },\n  ControlFlowTok: { bold: true },\n  InformationTok: { bold: true, italic: true },\n  WarningTok: { bold: true, italic: true },\n  AlertTok: { bold: true },\n  ErrorTok: { bold: true },\n};\n\nconst styles = context.document.getStyles();\n\nfor (const styleName of Object.keys(fixups)) {\n  const fixup = fixups[styleName];\n  const style = styles.getByNameOrNullObject(styleName);\n  style.load(\"isNullObject\");\n  await context.sync();\n  if (style.isNullObject) continue;\n\n  const font = style.font;\n  if (fixup.bold) {\n    font.bold = true;\n  }\n  if (fixup.italic) {\n    font.italic = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# fix(docx): fix OOXMLValidator error on KeywordTok output\n#\n# The <w:rPr> child-element order inside several syntax-highlighting\n# character styles in styles.xml doesn't match the sequence required by\n# wml.xsd (CT_RPr): <w:b/>/<w:i/> must come before <w:color/>. Re-assert\n# each style's bold/italic flag (it already has that value) so the host\n# re-serializes the style's <w:rPr> in canonical schema order.\n\n$d = $word.ActiveDocument\n\n# styleId -> which flag(s) to re-assert (only the flag(s) actually present\n# need touching; re-setting one forces the whole <w:rPr> to be rewritten in\n# schema order, fixing any other out-of-order siblings too).\n$fixups = @{\n    \"KeywordTok\"       = @{ Bold = $true;  Italic = $false }\n    \"ImportTok\"        = @{ Bold = $true;  Italic = $false }\n    \"CommentTok\"       = @{ Bold = $false; Italic = $true }\n    \"DocumentationTok\" = @{ Bold = $false; Italic = $true }\n    \"AnnotationTok\"    = @{ Bold = $true;  Italic = $true }\n    \"CommentVarTok\"    = @{ Bold = $true;  Italic = $true }\n    \"ControlFlowTok\"   = @{ Bold = $true;  Italic = $false }\n    \"InformationTok\"   = @{ Bold = $true;  Italic = $true }\n    \"WarningTok\"       = @{ Bold = $true;  Italic = $true }\n    \"AlertTok\"         = @{ Bold = $true;  Italic = $false }\n    \"ErrorTok\"         = @{ Bold = $true;  Italic = $false }\n}\n\nforeach ($styleId in $fixups.Keys) {\n    $fixup = $fixups[$styleId]\n    try {\n        $s = $d.Styles.Item($styleId)\n    } catch {\n        continue\n    }\n    if ($null -eq $s) {\n        continue\n    }\n    if ($fixup.Bold) {\n        $s.Font.Bold = -1\n    }\n    if ($fixup.Italic) {\n        $s.Font.Italic = -1\n    }\n}\n"}
